$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 11965
$ws1.Range("F6").Value = 361
$ws1.Range("F8").Value = 11859
$ws1.Range("F11").Value = 107
$ws1.Range("F12").Value = 74
$ws1.Range("F14").Value = 5874
$ws1.Range("F17").Value = 193

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 11965
$ws4.Range("F9").Value = 361
$ws4.Range("F11").Value = 11859
$ws4.Range("F14").Value = 107
$ws4.Range("F15").Value = 74
$ws4.Range("F18").Value = 5874
$ws4.Range("F21").Value = 193
